$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $value)
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws 'D2' '289.93'
Set-TextValue $ws 'E2' '-9.57%'
Set-TextValue $ws 'D3' '40.12'
Set-TextValue $ws 'E3' '-3.14%'
Set-TextValue $ws 'D4' '5.046'
Set-TextValue $ws 'E4' '-3.85%'
Set-TextValue $ws 'D5' '0.07298'
Set-TextValue $ws 'E5' '-5.76%'
Set-TextValue $ws 'D6' '4.278'
Set-TextValue $ws 'E6' '-1.40%'
Set-TextValue $ws 'D7' '1.553'
Set-TextValue $ws 'E7' '-10.70%'
Set-TextValue $ws 'D8' '0.9185'
Set-TextValue $ws 'D9' '0.1163'
Set-TextValue $ws 'E9' '-7.70%'
Set-TextValue $ws 'D10' '0.1733'
Set-TextValue $ws 'E10' '-7.05%'
Set-TextValue $ws 'D11' '0.08695'
Set-TextValue $ws 'E11' '-5.48%'
Set-TextValue $ws 'D12' '0.04168'
Set-TextValue $ws 'E12' '-0.03%'
Set-TextValue $ws 'D13' '0.1053'
Set-TextValue $ws 'E13' '0.12%'
Set-TextValue $ws 'D14' '0.001270'
Set-TextValue $ws 'E14' '-1.36%'
Set-TextValue $ws 'D15' '0.005778'
Set-TextValue $ws 'E15' '-1.05%'
Set-TextValue $ws 'D16' '3.392'
Set-TextValue $ws 'D19' '7.843'
Set-TextValue $ws 'E19' '-6.80%'
Set-TextValue $ws 'D20' '0.1351'
Set-TextValue $ws 'E20' '-0.25%'
Set-TextValue $ws 'D21' '0.2883'
Set-TextValue $ws 'E21' '5.55%'
Set-TextValue $ws 'D22' '0.03865'
Set-TextValue $ws 'E22' '-4.08%'
Set-TextValue $ws 'D23' '0.001269'
Set-TextValue $ws 'E23' '0.03%'
Set-TextValue $ws 'D24' '0.003882'
Set-TextValue $ws 'E24' '-5.95%'
Set-TextValue $ws 'D25' '0.0001281'
Set-TextValue $ws 'E25' '0.70%'
Set-TextValue $ws 'D26' '0.0003724'
Set-TextValue $ws 'D38' '0.02322'
Set-TextValue $ws 'E38' '-9.01%'
Set-TextValue $ws 'D39' '0.04963'
Set-TextValue $ws 'E39' '-7.18%'
Set-TextValue $ws 'D40' '0.006835'
Set-TextValue $ws 'E40' '243.08%'
Set-TextValue $ws 'D41' '0.007661'
Set-TextValue $ws 'E41' '-1.61%'
Set-TextValue $ws 'D42' '0.1274'
Set-TextValue $ws 'E42' '-3.33%'
Set-TextValue $ws 'D43' '0.007359'
Set-TextValue $ws 'E43' '4.69%'
Set-TextValue $ws 'D44' '0.007065'
Set-TextValue $ws 'E44' '-14.91%'
Set-TextValue $ws 'D45' '0.2888'
Set-TextValue $ws 'E45' '-16.64%'
Set-TextValue $ws 'D46' '0.00006413'
Set-TextValue $ws 'E46' '-4.34%'
Set-TextValue $ws 'D47' '0.00000000751'
Set-TextValue $ws 'E47' '-0.09%'
Set-TextValue $ws 'B48' 'BOLO'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws 'D48' '0.008538'
Set-TextValue $ws 'E48' '-95.70%'
Set-TextValue $ws 'B49' 'CoinbaseStockToken'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws 'D49' '0.004202'
Set-TextValue $ws 'E49' '-0.12%'
Set-TextValue $ws 'E50' '-0.09%'
Set-TextValue $ws 'E51' '-0.09%'
